$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.319.64"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.05%  '

$ws.Range('D3').Value = "'3.238.28"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.83%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = "'595.52"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.25%  '

$ws.Range('D6').Value = "'141.12"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.11%  '

$ws.Range('D8').Value = "'3.233.80"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.94%  '

$ws.Range('E9').Value = '  -1.82%  '

$ws.Range('E10').Value = '  -1.15%  '

$ws.Range('D11').Value = "'5.39"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.06%  '

$ws.Range('E12').Value = '  -0.45%  '

$ws.Range('E13').Value = '  -2.88%  '

$ws.Range('D14').Value = "'34.37"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.69%  '

$ws.Range('D15').Value = "'3.767.85"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.78%  '

$ws.Range('E16').Value = '  +0.07%  '

$ws.Range('D17').Value = "'3.232.76"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.41%  '

$ws.Range('D18').Value = "'63.344.55"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.02%  '

$ws.Range('D19').Value = "'6.80"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.93%  '

$ws.Range('D20').Value = "'473.98"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.04%  '

$ws.Range('D21').Value = "'14.20"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.35%  '

$ws.Range('D22').Value = "'0.731"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.62%  '

$ws.Range('D23').Value = "'7.94"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.48%  '

$ws.Range('D24').Value = "'83.57"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.50%  '

$ws.Range('E25').Value = '  -0.91%  '

$ws.Range('D27').Value = "'7.51"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.25%  '

$ws.Range('E28').Value = '  -0.97%  '

$ws.Range('D29').Value = "'8.13"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.01%  '

$ws.Range('D30').Value = "'2.13"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.73%  '

$ws.Range('D31').Value = "'27.45"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.91%  '

$ws.Range('E32').Value = '  -0.05%  '

$ws.Range('E33').Value = '  -4.43%  '

$ws.Range('D34').Value = "'2.54"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.32%  '

$ws.Range('E35').Value = '  -1.62%  '

$ws.Range('D36').Value = "'5.93"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.24%  '

$ws.Range('E37').Value = '  -0.11%  '

$ws.Range('D38').Value = "'0.0₃0709"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.63%  '

$ws.Range('E39').Value = '  -1.33%  '

$ws.Range('D40').Value = "'422.72"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.38%  '

$ws.Range('E41').Value = '  +0.25%  '

$ws.Range('D42').Value = "'2.75"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.33%  '

$ws.Range('D43').Value = "'2.968.40"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.19%  '

$ws.Range('D44').Value = "'0.110"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.67%  '

$ws.Range('E45').Value = '  +2.88%  '

$ws.Range('E46').Value = '  -1.17%  '

$ws.Range('E47').Value = '  +0.08%  '

$ws.Range('D48').Value = "'2.34"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.39%  '

$ws.Range('D49').Value = "'25.96"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.46%  '

$ws.Range('E50').Value = '  -0.53%  '

$ws.Range('D51').Value = "'121.06"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.22%  '
